$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 597 ("お金では買えないよ" post), shifting all
# subsequent rows up by one.
$ws.Rows.Item(597).Delete()
